# Student Remove Submission Function - DEVELOP
#
# Sample-data update accompanying the "remove submission" feature:
#   - Submission row for module 39596939 (student 59262392) now has an
#     actual marked report attached (was a stub row with reportId 0 and
#     status MARKED_2 / FINAL_YEAR) -> becomes report 18449474,
#     PENDING_MARKING, CAPSTONE_1.
#   - Submission row for module 37346231 (student 59262392) likewise now
#     has a real report attached (was reportId 0, OVERDUE, INVESTIGATION)
#     -> becomes report 63860114, PENDING_MARKING, CAPSTONE_2.
#   - Submission row for module 36887009 (student 59262392) reportType
#     corrected to REPORT.
#   - Report catalogue: row 2's type corrected to REPORT; rows 3 and 4
#     (InvestigateReport / FinalYearProject) had their name+path swapped
#     to line up with the correct report ids.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Submission sheet updates
# ---------------------------------------------------------------------
$submission = $wb.Worksheets.Item("Submission")

# Row 2 (id 2127241): reportType INVESTIGATION -> REPORT
$submission.Range("G2").Value = "REPORT"

# Row 11 (id 98294670): now has a real report attached
$submission.Range("B11").Value = "18449474"
$submission.Range("F11").Value = "PENDING_MARKING"
$submission.Range("G11").Value = "CAPSTONE_1"

# Row 14 (id 59626184): now has a real report attached
$submission.Range("B14").Value = "63860114"
$submission.Range("F14").Value = "PENDING_MARKING"
$submission.Range("G14").Value = "CAPSTONE_2"

# ---------------------------------------------------------------------
# Report sheet updates
# ---------------------------------------------------------------------
$report = $wb.Worksheets.Item("Report")

# Row 2 (id 49053257 / "Report"): reportType INVESTIGATION -> REPORT
$report.Range("D2").Value = "REPORT"

# Row 3 (id 26662640) and Row 4 (id 34685929): reportName/reportPath swapped
$report.Range("B3").Value = "FinalYearProject"
$report.Range("C3").Value = "src/main/resources/Data/SampleDataXlsx/FinalYearProject.pdf"
$report.Range("B4").Value = "InvestigateReport"
$report.Range("C4").Value = "src/main/resources/Data/SampleDataXlsx/InvestigateReport.pdf"

# ---------------------------------------------------------------------
# View-state: author was last working on the Intake sheet's D26, the
# Submission sheet (now the active tab, cell A8), the Consultation
# sheet (selection reset to A1) and the Report sheet (cell C21) when
# the file was saved.
# ---------------------------------------------------------------------
$intake = $wb.Worksheets.Item("Intake")
$intake.Activate()
$intake.Range("D26").Select()

$consultation = $wb.Worksheets.Item("Consultation")
$consultation.Activate()
$consultation.Range("A1").Select()

$report.Activate()
$report.Range("C21").Select()

$submission.Activate()
$submission.Range("A8").Select()
